{"js": "// Update the worksheet date header and the 25 multiplication problems\n// (5 problem-rows x 5 cells each) inside the single table, in document\n// order. Every problem cell's text changes; the table shape (rows/cols)\n// itself is unchanged, so we only need to overwrite run text positionally.\n\n// 1) Update the date/weekday heading line.\nconst dateResults = context.document.body.search(\"2024-06-16 Sunday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2024-06-17 Monday\", Word.InsertLocation.replace);\n}\n\n// 2) Update the multiplication problems table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// New values for each of the 5 \"problem\" rows (in the order they appear\n// in the document). Rows that hold no problem text (the blank work-space\n// rows in between) are left untouched.\nconst newRowValues = [\n  [\"628\u00d78=\", \"814\u00d75=\", \"369\u00d76=\", \"574\u00d77=\", \"680\u00d72=\"],\n  [\"420\u00d79=\", \"762\u00d73=\", \"518\u00d73=\", \"434\u00d78=\", \"924\u00d76=\"],\n  [\"281\u00d79=\", \"293\u00d78=\", \"167\u00d73=\", \"700\u00d72=\", \"381\u00d78=\"],\n  [\"998\u00d72=\", \"815\u00d72=\", \"310\u00d73=\", \"804\u00d73=\", \"115\u00d74=\"],\n  [\"778\u00d73=\", \"547\u00d72=\", \"914\u00d76=\", \"525\u00d78=\", \"628\u00d79=\"],\n];\n\n// Load every row's cells up front.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// For every cell, also load its text so we can detect which rows are\n// \"problem\" rows (non-blank first cell) vs. blank spacer rows.\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.load(\"value\");\n  }\n}\nawait context.sync();\n\nlet problemRowIndex = 0;\nfor (const row of rows.items) {\n  const cells = row.cells.items;\n  const isProblemRow = cells.length > 0 && cells[0].value.trim() !== \"\";\n  if (!isProblemRow) {\n    continue;\n  }\n  const values = newRowValues[problemRowIndex];\n  if (values) {\n    for (let c = 0; c < cells.length && c < values.length; c++) {\n      cells[c].value = values[c];\n    }\n  }\n  problemRowIndex++;\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date header and the 25 multiplication problems\n# (5 problem-rows x 5 cells each) inside the single table, in document\n# order. Every problem cell's text changes; the table shape (rows/cols)\n# itself is unchanged, so we only need to overwrite cell text positionally.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date/weekday heading line.\n$find = $d.Content.Find\n$find.Text = \"2024-06-16 Sunday\"\n$find.Replacement.Text = \"2024-06-17 Monday\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2) Update the multiplication problems table.\n$t = $d.Tables.Item(1)\n\n# New values for each of the 5 \"problem\" rows (in the order they appear\n# in the document). Rows that hold no problem text (the blank work-space\n# rows in between) are left untouched.\n$newRowValues = @(\n  @(\"628\u00d78=\", \"814\u00d75=\", \"369\u00d76=\", \"574\u00d77=\", \"680\u00d72=\"),\n  @(\"420\u00d79=\", \"762\u00d73=\", \"518\u00d73=\", \"434\u00d78=\", \"924\u00d76=\"),\n  @(\"281\u00d79=\", \"293\u00d78=\", \"167\u00d73=\", \"700\u00d72=\", \"381\u00d78=\"),\n  @(\"998\u00d72=\", \"815\u00d72=\", \"310\u00d73=\", \"804\u00d73=\", \"115\u00d74=\"),\n  @(\"778\u00d73=\", \"547\u00d72=\", \"914\u00d76=\", \"525\u00d78=\", \"628\u00d79=\")\n)\n\n$problemRowIndex = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  $firstCell = $t.Cell($r, 1)\n  # Word always reports at least a paragraph mark (chr 13) plus a\n  # trailing cell mark (chr 7) for a table cell's Range.Text, even when\n  # the cell is visually empty, so strip those control characters (not\n  # just whitespace) before testing for blankness.\n  $firstCellText = ($firstCell.Range.Text -replace \"[\\r\\a]\", \"\")\n  if ($firstCellText -eq \"\") {\n    continue\n  }\n  $values = $newRowValues[$problemRowIndex]\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $values[$c - 1]\n  }\n  $problemRowIndex++\n}\n"}
